# feat: add 2022-Q1 data
#
# - insert a new "2022-Q1" sheet (holdings detail) right before the "总计"
#   (total) summary sheet, copying the layout/style of the most recent
#   quarter sheet ("2021-Q4")
# - update the "总计" summary sheet with a new leading row for 2022-Q1 and
#   renumber/relabel the existing rows

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")
$template = $wb.Worksheets.Item("2021-Q4")

# 1) Clone the latest quarter sheet (keeps fonts/borders/number formats
#    identical to the other quarter tabs) and drop the clone right before
#    the total sheet, then rename it. (Re-resolve the sheet right after the
#    copy instead of reusing $totalSheet - once a sheet is inserted, stale
#    worksheet handles keep resolving by their old numeric position.)
$template.Copy($totalSheet)
$newQuarter = $wb.Worksheets.Item($template.Index + 1)
$newQuarter.Name = "2022-Q1"

# 2) Overwrite the single holding row with the 2022-Q1 figures. The
#    numeric-looking columns are stored as text (matching every other
#    quarter sheet), so force a text number format before assigning.
$newQuarter.Range("B2").NumberFormat = "@"
$newQuarter.Range("B2").Value = "162416"
$newQuarter.Range("C2").NumberFormat = "@"
$newQuarter.Range("C2").Value = "华宝港股通恒生香港35指数(LOF)"
$newQuarter.Range("D2").NumberFormat = "@"
$newQuarter.Range("D2").Value = "0.21"
$newQuarter.Range("E2").NumberFormat = "@"
$newQuarter.Range("E2").Value = "94.50"
$newQuarter.Range("F2").NumberFormat = "@"
$newQuarter.Range("F2").Value = "4.45"
$newQuarter.Range("G2").NumberFormat = "@"
$newQuarter.Range("G2").Value = "0.0093"
$newQuarter.Range("H2").Value = 7

# 3) Restore the originally-active tab (the first sheet) since inserting /
#    copying sheets shifts the active selection to the new tab.
$wb.Worksheets.Item(1).Activate()

# 4) Update the "总计" sheet: insert a fresh row for 2022-Q1 above the
#    existing rows, then renumber the index column / relabel the quarters
#    that shifted down by one row. Re-resolve the sheet by name again since
#    the worksheet collection changed above.
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows(2).Insert()
$totalSheet.Rows(2).ClearFormats()

$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.01

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q4"
$totalSheet.Range("C3").Value = 1
$totalSheet.Range("D3").Value = 0.01

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2021-Q3"
$totalSheet.Range("C4").Value = 1
$totalSheet.Range("D4").Value = 0.01

$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = "2021-Q2"
$totalSheet.Range("C5").Value = 1
$totalSheet.Range("D5").Value = 0.01

$totalSheet.Range("A6").Value = 4
$totalSheet.Range("B6").Value = "2021-Q1"
$totalSheet.Range("C6").Value = 2
$totalSheet.Range("D6").Value = 0.02

$totalSheet.Range("A7").Value = 5
$totalSheet.Range("B7").Value = "2020-Q4"
$totalSheet.Range("C7").Value = 1
$totalSheet.Range("D7").Value = 0.01
